$wb = $excel.ActiveWorkbook

# --- Sheet: Recommandations ---
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws1.Range("A2").Value = 'BRVM-PRINCIPAL     (**)'
$ws1.Range("D2").Value = 950.76
$ws1.Range("E2").Value = 239.91
$ws1.Range("A3").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 719.84
$ws1.Range("E3").Value = 181.57
$ws1.Range("A4").Value = 'BRVM - CONSOMMATION DE BASE     (**)'
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 706.3099999999999
$ws1.Range("E4").Value = 237.88
$ws1.Range("A5").Value = 'BRVM - INDUSTRIELS'
$ws1.Range("D5").Value = 659.95
$ws1.Range("E5").Value = 171.31
$ws1.Range("A6").Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Range("D6").Value = 625.24
$ws1.Range("E6").Value = 157.36
$ws1.Range("A7").Value = 'BRVM-PRESTIGE'
$ws1.Range("D7").Value = 603.58
$ws1.Range("E7").Value = 152.04
$ws1.Range("A8").Value = 'BRVM - ENERGIE'
$ws1.Range("D8").Value = 479.02
$ws1.Range("E8").Value = 121.85
$ws1.Range("A9").Value = 'BRVM - SERVICES PUBLICS'
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 470.93
$ws1.Range("E9").Value = 118.28
$ws1.Range("A10").Value = 'BRVM – COMPOSITE TOTAL RETURN     (**)'
$ws1.Range("C10").Value = 3
$ws1.Range("D10").Value = 424.58
$ws1.Range("E10").Value = 142.33
$ws1.Range("A11").Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 391.76
$ws1.Range("E11").Value = 98.31
$ws1.Range("D14").Value = 26.51
$ws1.Range("E14").Value = 7.31
$ws1.Range("A16").Value = 'SICABLE CI (CABC)'
$ws1.Range("B16").Value = 2
$ws1.Range("D16").Value = 14.17
$ws1.Range("E16").Value = 7.5
$ws1.Range("A17").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("D17").Value = 7.27
$ws1.Range("E17").Value = 7.27
$ws1.Range("A18").Value = 'SICOR CI (SICC)'
$ws1.Range("B18").Value = 2
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = 6.37
$ws1.Range("E18").Value = 7.47
$ws1.Range("G18").Value = '👀 À surveiller'
$ws1.Range("A19").Value = 'BERNABE CI (BNBC)'
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 3.54
$ws1.Range("E19").Value = -3.25
$ws1.Range("G19").Value = '👀 À surveiller'
$ws1.Range("A20").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("D20").Value = 2.71
$ws1.Range("E20").Value = -3.1
$ws1.Range("A21").Value = 'SODE CI (SDCC)'
$ws1.Range("D21").Value = 2.59
$ws1.Range("E21").Value = 7.39
$ws1.Range("A22").Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Range("B22").Value = 1
$ws1.Range("D22").Value = 2.16
$ws1.Range("E22").Value = -2.04
$ws1.Range("A23").Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Range("D23").Value = 1.81
$ws1.Range("E23").Value = -2.09
$ws1.Range("A25").Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Range("B25").Value = 0
$ws1.Range("D25").Value = -1.02
$ws1.Range("E25").Value = -1.02
$ws1.Range("G25").Value = '➖ Neutre'
$ws1.Range("A26").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("D26").Value = -1.7
$ws1.Range("E26").Value = -1.7
$ws1.Range("A27").Value = 'CIE CI (CIEC)'
$ws1.Range("D27").Value = -1.86
$ws1.Range("E27").Value = -1.86
$ws1.Range("A28").Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Range("D28").Value = -2
$ws1.Range("E28").Value = -2
$ws1.Range("A29").Value = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$ws1.Range("B29").Value = 1
$ws1.Range("D29").Value = -2.21
$ws1.Range("E29").Value = 2.44
$ws1.Range("G29").Value = '👀 À surveiller'
$ws1.Range("A30").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("D30").Value = -2.21
$ws1.Range("E30").Value = -2.21
$ws1.Range("A31").Value = 'SAPH CI (SPHC)'
$ws1.Range("B31").Value = 0
$ws1.Range("D31").Value = -2.4
$ws1.Range("E31").Value = -2.4
$ws1.Range("G31").Value = '➖ Neutre'
$ws1.Range("A32").Value = 'BANK OF AFRICA ML (BOAM)'
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = -2.41
$ws1.Range("E32").Value = -2.41
$ws1.Range("G32").Value = '➖ Neutre'
$ws1.Range("A33").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("D33").Value = -3.34
$ws1.Range("E33").Value = -3.34
$ws1.Range("A34").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B34").Value = 1
$ws1.Range("D34").Value = -3.51
$ws1.Range("E34").Value = 2.71
$ws1.Range("G34").Value = '👀 À surveiller'
$ws1.Range("A35").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -4.17
$ws1.Range("E35").Value = -4.17

# --- Sheet: Top_YTD ---
$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Range("A2").Value = 'BRVM-PRINCIPAL     (**)'
$ws2.Range("B2").Value = 12903.17
$ws2.Range("A3").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws2.Range("B3").Value = 6042.86
$ws2.Range("A4").Value = 'BRVM - INDUSTRIELS'
$ws2.Range("B4").Value = 4828.1
$ws2.Range("A5").Value = 'BRVM - SERVICES FINANCIERS'
$ws2.Range("B5").Value = 4215.72
$ws2.Range("A6").Value = 'BRVM-PRESTIGE'
$ws2.Range("B6").Value = 3862.42
$ws2.Range("A7").Value = 'BRVM - CONSOMMATION DE BASE     (**)'
$ws2.Range("B7").Value = 3673.73
$ws2.Range("A8").Value = 'BRVM - ENERGIE'
$ws2.Range("B8").Value = 2231.92
$ws2.Range("A9").Value = 'BRVM - SERVICES PUBLICS'
$ws2.Range("B9").Value = 2147.33
$ws2.Range("A10").Value = 'BRVM - TELECOMMUNICATIONS'
$ws2.Range("B10").Value = 1435.08
$ws2.Range("B11").Value = 1308.93
